$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Donor cell D47 keeps its original default (unstyled) format throughout;
# used to restore number formatting on D-cells that would otherwise be
# auto-coerced to numeric type by Excel when the new text looks like a plain number.

$ws.Range("D2").Value = "68.829.90"
$ws.Range("E2").Value = "  +2.15%  "

$ws.Range("D3").Value = "3.731.01"
$ws.Range("E3").Value = "  -1.06%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.86"
$ws.Range("D47").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("E5").Value = "  +1.08%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.33"
$ws.Range("D47").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("E6").Value = "  -4.84%  "

$ws.Range("D7").Value = "3.731.09"
$ws.Range("E7").Value = "  -0.94%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.533"
$ws.Range("D47").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("E9").Value = "  +2.89%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.163"
$ws.Range("D47").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("E10").Value = "  +3.16%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.35"
$ws.Range("D47").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("E11").Value = "  +2.98%  "

$ws.Range("E12").Value = "  -0.31%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.00"
$ws.Range("D47").Copy()
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("E13").Value = "  -0.08%  "

$ws.Range("E14").Value = "  +0.05%  "

$ws.Range("D15").Value = "4.356.68"
$ws.Range("E15").Value = "  -1.57%  "

$ws.Range("D16").Value = "3.735.15"
$ws.Range("E16").Value = "  -2.19%  "

$ws.Range("D17").Value = "68.816.47"
$ws.Range("E17").Value = "  +2.05%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.24"
$ws.Range("D47").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E18").Value = "  +1.71%  "

$ws.Range("E19").Value = "  +0.54%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.24"
$ws.Range("D47").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").Value = "  +5.96%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "496.72"
$ws.Range("D47").Copy()
$ws.Range("D21").PasteSpecial(-4122)
$ws.Range("E21").Value = "  +1.99%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.17"
$ws.Range("D47").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = "  +12.52%  "

$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.14"
$ws.Range("D47").Copy()
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("E24").Value = "  +1.61%  "

$ws.Range("E25").Value = "  -1.32%  "

$ws.Range("E26").Value = "  -6.00%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.35"
$ws.Range("D47").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = "  +1.49%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.11"
$ws.Range("D47").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = "  -0.52%  "

$ws.Range("E29").Value = "  -0.13%  "

$ws.Range("E30").Value = "  +0.31%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.45"
$ws.Range("D47").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").Value = "  +1.03%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.97"
$ws.Range("D47").Copy()
$ws.Range("D32").PasteSpecial(-4122)
$ws.Range("E32").Value = "  +3.33%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.66"
$ws.Range("D47").Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("E33").Value = "  -3.27%  "

$ws.Range("D34").Value = "3.882.48"
$ws.Range("E34").Value = "  -1.36%  "

$ws.Range("E35").Value = "  -0.17%  "

$ws.Range("D36").Value = "3.664.01"
$ws.Range("E36").Value = "  -1.70%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D47").Copy()
$ws.Range("D37").PasteSpecial(-4122)
$ws.Range("E37").Value = "  -0.12%  "

$ws.Range("E38").Value = "  +1.14%  "

$ws.Range("E39").Value = "  +1.34%  "

$ws.Range("E40").Value = "  -1.94%  "

$ws.Range("E41").Value = "  -0.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "434.73"
$ws.Range("D47").Copy()
$ws.Range("D42").PasteSpecial(-4122)
$ws.Range("E42").Value = "  -3.89%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.88"
$ws.Range("D47").Copy()
$ws.Range("D43").PasteSpecial(-4122)
$ws.Range("E43").Value = "  -0.26%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.98"
$ws.Range("D47").Copy()
$ws.Range("D44").PasteSpecial(-4122)
$ws.Range("E44").Value = "  -0.44%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.85"
$ws.Range("D47").Copy()
$ws.Range("D45").PasteSpecial(-4122)
$ws.Range("E45").Value = "  +0.25%  "

$ws.Range("E46").Value = "  +1.29%  "

$ws.Range("E48").Value = "  -1.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "142.06"
$ws.Range("D47").Copy()
$ws.Range("D49").PasteSpecial(-4122)
$ws.Range("E49").Value = "  +0.42%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0352"
$ws.Range("D47").Copy()
$ws.Range("D50").PasteSpecial(-4122)
$ws.Range("E50").Value = "  +0.98%  "

$ws.Range("D51").Value = "2.739.18"
$ws.Range("E51").Value = "  -2.71%  "

$excel.CutCopyMode = 0

